# ToDo.xlsx: add "Conexion a interfaz" sheet + tweak selections
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Move the selection on "Entregables" (sheet 1) from C14 to C23
# ---------------------------------------------------------------------
$wsEntregables = $wb.Worksheets.Item("Entregables")
$wsEntregables.Activate()
$wsEntregables.Range("C23").Select()

# ---------------------------------------------------------------------
# 2) Add the new worksheet "Conexion a interfaz" as the last (4th) tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Name = "Conexion a interfaz"

# ---------------------------------------------------------------------
# 3) Populate the new sheet
# ---------------------------------------------------------------------
$rows = @(
    @("Camino mas corto de min desde una habitacion hasta la salida", $true),
    @("Camino mas corto entre dos habitaciones ", $false),
    @("Transmitir mensaje de cierre", $false),
    @("Añadir una habitacion", $true),
    @("Eliminar habitacion", $true),
    @("Registrar tesoros encontrados", $false),
    @("Visualizar los tesoros", $false),
    @("Añadir pasillo", $true),
    @("Eliminar Pasillo", $true)
)

$r = 1
foreach ($row in $rows) {
    $wsNew.Range("A$r").Value = $row[0]
    if ($row[1]) {
        $wsNew.Range("B$r").Value = "x"
        $wsNew.Range("B$r").HorizontalAlignment = -4108
    }
    $r = $r + 1
}

$wsNew.Columns.Item(1).AutoFit() | Out-Null

# Conditional formatting: highlight "x" marks in column B with the
# built-in "Green Fill with Dark Green Text" style.
# NOTE: OLE_COLOR values are packed as 0x00BBGGRR (BGR), not RGB.
$fc = $wsNew.Range("B1:B1048576").FormatConditions
$rule = $fc.Add(8, 3, '"x"')
$rule.Font.Color = 0x00006100      # R=0x00 G=0x61 B=0x00
$rule.Interior.Color = 0x00CEEFC6  # R=0xC6 G=0xEF B=0xCE

# Leave the cursor where the author left it
$wsNew.Range("D4").Select()
